$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.360.91'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.355.85'
$ws.Range("E3").Value = '  +2.93%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.49'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.97'
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("E9").Value = '  +3.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.14'
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.86'
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0814'
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.00'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").Value = '2.721.04'
$ws.Range("E15").Value = '  +3.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.72'
$ws.Range("E16").Value = '  +6.22%  '
$ws.Range("D17").Value = '2.357.30'
$ws.Range("E17").Value = '  +3.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.813'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("D19").Value = '43.348.89'
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("E20").Value = '  -4.02%  '
$ws.Range("D21").Value = '0.0₃0929'
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("E22").Value = '  +3.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.32'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '242.59'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.06'
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.95'
$ws.Range("E28").Value = '  +8.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.31'
$ws.Range("E29").Value = '  +7.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.65'
$ws.Range("E30").Value = '  -5.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.61'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.68'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.31'
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("E36").Value = '  +6.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.11'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.64'
$ws.Range("E39").Value = '  +10.70%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("E40").Value = '  +5.78%  '
$ws.Range("E41").Value = '  +1.28%  '
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.44'
$ws.Range("E43").Value = '  +5.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.76'
$ws.Range("E44").Value = '  +1.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0292'
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("D46").Value = '1.992.28'
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.10'
$ws.Range("E47").Value = '  +2.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.45'
$ws.Range("E48").Value = '  +6.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '58.49'
$ws.Range("E49").Value = '  +6.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.92'
$ws.Range("E50").Value = '  -3.01%  '
$ws.Range("E51").Value = '  +2.93%  '
